# UMS_Data.xlsx - "Subjects" sheet update.
#
# Row 9 previously held a leftover test entry ("FACTORIO204" /
# "resource management"); it now holds the real subject that belongs
# there ("ENG1210" / "Intro to Coding"). A new row 10 adds the
# "FORTNITE101" / "fortnite" subject that was missing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Subjects")

$ws.Range("A9").Value = "ENG1210"
$ws.Range("B9").Value = "Intro to Coding"

$ws.Range("A10").Value = "FORTNITE101"
$ws.Range("B10").Value = "fortnite"
